$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Original long query (was in B2, originally the DatasetsTab query) - now moves to B3 (ProjectsTab)
$origQuery = @"
SELECT DISTINCT
    ds.dataset_title AS "Title",
    ds.dataset_source_id AS "Source ID", 
    ds.primary_disease AS "Primary Disease",
    CAST(ds.participant_count AS INT) AS "Participants Count",
    CAST(ds.sample_count AS INT) AS "Sample Count",
    CASE 
        WHEN LENGTH(REPLACE(TRIM(ds.description), '  ', ' ')) > 500 
        THEN SUBSTR(REPLACE(TRIM(ds.description), '  ', ' '), 1, 500) || ' ...'
        ELSE REPLACE(TRIM(ds.description), '  ', ' ')
    END AS "Description"
FROM df_cedcd ds
ORDER BY ds.dataset_title ASC;
"@

# New, shorter query for DatasetsTab (B2)
$newQuery = @"
SELECT DISTINCT
    REPLACE(ds.dataset_title, '  ', ' ') AS "Title", 
    ds.dataset_source_id AS "Source ID", 
    ds.primary_disease AS "Primary Disease",
    -- CAST(ds.participant_count AS INT) AS "Participants Count",
    CAST(ds.sample_count AS INT) AS "Sample Count"
FROM df_cedcd ds
ORDER BY ds.dataset_title ASC;
"@

$ws.Range("B2").Value = $newQuery
$ws.Range("B3").Value = $origQuery

$ws.Range("C2").Select()
